$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Differences in Tax Wedges")

# Row 9
$ws.Range("B9").Value = 2023
$ws.Range("C9").Value = 2024
$ws.Range("D9").Value = 2025
$ws.Range("E9").Value = 2026
$ws.Range("F9").Value = 2027
$ws.Range("G9").Value = 2028
$ws.Range("H9").Value = 2029
$ws.Range("I9").Value = 2030
$ws.Range("J9").Value = 2031
$ws.Range("K9").Value = 2032
$ws.Range("L9").Value = 2033

# Row 11
$ws.Range("B11").Value = 1.04
$ws.Range("C11").Value = 1.08
$ws.Range("D11").Value = 1.11
$ws.Range("E11").Value = 1.21
$ws.Range("F11").Value = 1.28
$ws.Range("G11").Value = 1.27
$ws.Range("H11").Value = 1.24
$ws.Range("I11").Value = 1.25
$ws.Range("J11").Value = 1.25
$ws.Range("K11").Value = 1.26
$ws.Range("L11").Value = 1.26

# Row 14
$ws.Range("B14").Value = 1.26
$ws.Range("C14").Value = 1.3
$ws.Range("D14").Value = 1.34
$ws.Range("E14").Value = 1.51
$ws.Range("F14").Value = 1.58
$ws.Range("G14").Value = 1.58
$ws.Range("H14").Value = 1.55
$ws.Range("I14").Value = 1.55
$ws.Range("J14").Value = 1.55
$ws.Range("K14").Value = 1.55
$ws.Range("L14").Value = 1.55

# Row 15
$ws.Range("B15").Value = 0.59
$ws.Range("C15").Value = 0.63
$ws.Range("D15").Value = 0.63
$ws.Range("E15").Value = 0.6
$ws.Range("F15").Value = 0.64
$ws.Range("G15").Value = 0.64
$ws.Range("H15").Value = 0.61
$ws.Range("I15").Value = 0.62
$ws.Range("J15").Value = 0.63
$ws.Range("K15").Value = 0.65
$ws.Range("L15").Value = 0.66

# Row 17
$ws.Range("B17").Value = 0.67
$ws.Range("C17").Value = 0.68
$ws.Range("D17").Value = 0.72
$ws.Range("E17").Value = 0.91
$ws.Range("F17").Value = 0.94
$ws.Range("G17").Value = 0.94
$ws.Range("H17").Value = 0.94
$ws.Range("I17").Value = 0.93
$ws.Range("J17").Value = 0.92
$ws.Range("K17").Value = 0.9
$ws.Range("L17").Value = 0.89

# Row 20
$ws.Range("B20").Value = 1.17
$ws.Range("C20").Value = 1.2
$ws.Range("D20").Value = 1.23
$ws.Range("E20").Value = 1.27
$ws.Range("F20").Value = 1.32
$ws.Range("G20").Value = 1.33
$ws.Range("H20").Value = 1.3
$ws.Range("I20").Value = 1.31
$ws.Range("J20").Value = 1.31
$ws.Range("K20").Value = 1.32
$ws.Range("L20").Value = 1.32

# Row 21
$ws.Range("B21").Value = 0.81
$ws.Range("C21").Value = 0.86
$ws.Range("D21").Value = 0.89
$ws.Range("E21").Value = 1.11
$ws.Range("F21").Value = 1.19
$ws.Range("G21").Value = 1.18
$ws.Range("H21").Value = 1.12
$ws.Range("I21").Value = 1.13
$ws.Range("J21").Value = 1.13
$ws.Range("K21").Value = 1.15
$ws.Range("L21").Value = 1.15

# Row 23
$ws.Range("B23").Value = 0.36
$ws.Range("C23").Value = 0.34
$ws.Range("D23").Value = 0.34
$ws.Range("E23").Value = 0.16
$ws.Range("F23").Value = 0.14
$ws.Range("G23").Value = 0.15
$ws.Range("H23").Value = 0.18
$ws.Range("I23").Value = 0.17
$ws.Range("J23").Value = 0.18
$ws.Range("K23").Value = 0.17
$ws.Range("L23").Value = 0.17

# Row 25
$ws.Range("B25").Value = 0.8
$ws.Range("C25").Value = 0.75
$ws.Range("D25").Value = 0.71
$ws.Range("E25").Value = 0.7
$ws.Range("F25").Value = 0.66
$ws.Range("G25").Value = 0.66
$ws.Range("H25").Value = 0.65
$ws.Range("I25").Value = 0.65
$ws.Range("J25").Value = 0.65
$ws.Range("K25").Value = 0.65
$ws.Range("L25").Value = 0.65

# Row 26
$ws.Range("B26").Value = 0.08
$ws.Range("C26").Value = 0.07
$ws.Range("D26").Value = 0.06
$ws.Range("E26").Value = 0.06
$ws.Range("F26").Value = 0.02
$ws.Range("G26").Value = 0.02
$ws.Range("H26").Value = 0.02
$ws.Range("I26").Value = 0.02
$ws.Range("J26").Value = 0.02
$ws.Range("K26").Value = 0.02
$ws.Range("L26").Value = 0.02

# Row 28
$ws.Range("B28").Value = 0.05
$ws.Range("C28").Value = 0.02
$ws.Range("D28").Value = -0.02
$ws.Range("E28").Value = -0.5
$ws.Range("F28").Value = -0.51
$ws.Range("G28").Value = -0.5
$ws.Range("H28").Value = -0.49
$ws.Range("I28").Value = -0.49
$ws.Range("J28").Value = -0.49
$ws.Range("K28").Value = -0.49
$ws.Range("L28").Value = -0.49

# Row 29
$ws.Range("B29").Value = -1
$ws.Range("C29").Value = -1.07
$ws.Range("D29").Value = -1.13
$ws.Range("E29").Value = -1.71
$ws.Range("F29").Value = -1.78
$ws.Range("G29").Value = -1.78
$ws.Range("H29").Value = -1.73
$ws.Range("I29").Value = -1.73
$ws.Range("J29").Value = -1.73
$ws.Range("K29").Value = -1.75
$ws.Range("L29").Value = -1.75

# Row 32
$ws.Range("B32").Value = -0.02
$ws.Range("C32").Value = -0.02
$ws.Range("D32").Value = -0.02
$ws.Range("E32").Value = -0.25
$ws.Range("F32").Value = -0.26
$ws.Range("G32").Value = -0.25
$ws.Range("H32").Value = -0.24
$ws.Range("I32").Value = -0.24
$ws.Range("J32").Value = -0.24
$ws.Range("K32").Value = -0.25
$ws.Range("L32").Value = -0.24

# Row 33
$ws.Range("B33").Value = 0.16
$ws.Range("C33").Value = 0.08
$ws.Range("D33").Value = -0.01
$ws.Range("E33").Value = -0.93
$ws.Range("F33").Value = -0.94
$ws.Range("G33").Value = -0.94
$ws.Range("H33").Value = -0.91
$ws.Range("I33").Value = -0.91
$ws.Range("J33").Value = -0.91
$ws.Range("K33").Value = -0.92
$ws.Range("L33").Value = -0.92

# Row 35
$ws.Range("B35").Value = -0.18
$ws.Range("C35").Value = -0.09
$ws.Range("D35").Value = -0.01
$ws.Range("E35").Value = 0.68
$ws.Range("F35").Value = 0.69
$ws.Range("G35").Value = 0.69
$ws.Range("H35").Value = 0.66
$ws.Range("I35").Value = 0.67
$ws.Range("J35").Value = 0.67
$ws.Range("K35").Value = 0.67
$ws.Range("L35").Value = 0.67

# Row 38
$ws.Range("B38").Value = 0.05
$ws.Range("C38").Value = 0.02
$ws.Range("D38").Value = -0.02
$ws.Range("E38").Value = -0.5
$ws.Range("F38").Value = -0.51
$ws.Range("G38").Value = -0.5
$ws.Range("H38").Value = -0.49
$ws.Range("I38").Value = -0.49
$ws.Range("J38").Value = -0.49
$ws.Range("K38").Value = -0.49
$ws.Range("L38").Value = -0.49

# Row 39
$ws.Range("B39").Value = 1.42
$ws.Range("C39").Value = 1.41
$ws.Range("D39").Value = 1.39
$ws.Range("E39").Value = 1.45
$ws.Range("F39").Value = 1.45
$ws.Range("G39").Value = 1.44
$ws.Range("H39").Value = 1.41
$ws.Range("I39").Value = 1.41
$ws.Range("J39").Value = 1.41
$ws.Range("K39").Value = 1.43
$ws.Range("L39").Value = 1.43

# Row 41
$ws.Range("B41").Value = -1.37
$ws.Range("C41").Value = -1.39
$ws.Range("D41").Value = -1.4
$ws.Range("E41").Value = -1.95
$ws.Range("F41").Value = -1.95
$ws.Range("G41").Value = -1.95
$ws.Range("H41").Value = -1.9
$ws.Range("I41").Value = -1.9
$ws.Range("J41").Value = -1.9
$ws.Range("K41").Value = -1.92
$ws.Range("L41").Value = -1.92

